$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to stay as text,
# matching the source data (European-style / truncated numeric strings).
# Format them as Text before assigning so Excel does not coerce them to numbers.
$textCells = @("D5","D6","D8","D11","D12","D14","D16","D18","D19","D20","D21","D23","D24","D27","D29","D30","D32","D33","D34","D35","D36","D38","D39","D40","D41","D42","D43","D45","D48","D49","D50","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

# Apply the updated values
$ws.Range("D2").Value = "62.031.69"
$ws.Range("E2").Value = "  -2.10%  "
$ws.Range("D3").Value = "2.587.19"
$ws.Range("E3").Value = "  -4.74%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Value = "553.85"
$ws.Range("E5").Value = "  -1.26%  "
$ws.Range("D6").Value = "154.57"
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "0.593"
$ws.Range("E8").Value = "  +0.43%  "
$ws.Range("E9").Value = "  -2.54%  "
$ws.Range("E10").Value = "  -3.02%  "
$ws.Range("D11").Value = "5.48"
$ws.Range("E11").Value = "  -1.69%  "
$ws.Range("D12").Value = "0.365"
$ws.Range("E12").Value = "  -1.60%  "
$ws.Range("D13").Value = "3.047.61"
$ws.Range("E13").Value = "  -4.51%  "
$ws.Range("D14").Value = "25.55"
$ws.Range("E14").Value = "  -3.26%  "
$ws.Range("D15").Value = "61.913.10"
$ws.Range("E15").Value = "  -2.09%  "
$ws.Range("D16").Value = "0.0000143"
$ws.Range("E16").Value = "  -2.50%  "
$ws.Range("D17").Value = "2.591.65"
$ws.Range("E17").Value = "  -4.51%  "
$ws.Range("D18").Value = "11.65"
$ws.Range("E18").Value = "  -4.22%  "
$ws.Range("D19").Value = "4.55"
$ws.Range("E19").Value = "  -2.56%  "
$ws.Range("D20").Value = "339.00"
$ws.Range("E20").Value = "  -3.40%  "
$ws.Range("D21").Value = "6.04"
$ws.Range("E21").Value = "  -6.52%  "
$ws.Range("E22").Value = "  -0.21%  "
$ws.Range("D23").Value = "0.499"
$ws.Range("E23").Value = "  -2.61%  "
$ws.Range("D24").Value = "62.68"
$ws.Range("E24").Value = "  -2.61%  "
$ws.Range("E25").Value = "  -0.51%  "
$ws.Range("E26").Value = "  +0.00%  "
$ws.Range("D27").Value = "8.07"
$ws.Range("E27").Value = "  -1.47%  "
$ws.Range("D28").Value = "0.0₃0837"
$ws.Range("E28").Value = "  -6.03%  "
$ws.Range("D29").Value = "1.92"
$ws.Range("E29").Value = "  -1.37%  "
$ws.Range("D30").Value = "7.10"
$ws.Range("E30").Value = "  -0.98%  "
$ws.Range("E31").Value = "  -3.72%  "
$ws.Range("D32").Value = "160.06"
$ws.Range("E32").Value = "  -3.59%  "
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.02%  "
$ws.Range("D34").Value = "4.71"
$ws.Range("E34").Value = "  -2.30%  "
$ws.Range("D35").Value = "19.24"
$ws.Range("E35").Value = "  -3.03%  "
$ws.Range("D36").Value = "1.42"
$ws.Range("E36").Value = "  -3.92%  "
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").Value = "340.69"
$ws.Range("E38").Value = "  -1.27%  "
$ws.Range("D39").Value = "6.00"
$ws.Range("E39").Value = "  -1.97%  "
$ws.Range("D40").Value = "0.895"
$ws.Range("E40").Value = "  -7.18%  "
$ws.Range("D41").Value = "3.92"
$ws.Range("E41").Value = "  -3.07%  "
$ws.Range("D42").Value = "37.51"
$ws.Range("E42").Value = "  -2.43%  "
$ws.Range("D43").Value = "20.56"
$ws.Range("E43").Value = "  -3.80%  "
$ws.Range("E44").Value = "  +0.02%  "
$ws.Range("D45").Value = "0.608"
$ws.Range("E45").Value = "  -2.55%  "
$ws.Range("D46").Value = "2.136.11"
$ws.Range("E46").Value = "  +1.50%  "
$ws.Range("E47").Value = "  -4.99%  "
$ws.Range("B48").Value = "WhiteBITCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D48").Value = "10.94"
$ws.Range("E48").Value = "  -1.04%  "
$ws.Range("B49").Value = "Hedera"
$ws.Range("C49").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D49").Value = "0.0548"
$ws.Range("E49").Value = "  -4.38%  "
$ws.Range("D50").Value = "0.0965"
$ws.Range("E50").Value = "  -1.91%  "
$ws.Range("D51").Value = "0.0240"
$ws.Range("E51").Value = "  -2.46%  "
